# Add newly-translated skill/cube entries (Korean source text in column A,
# English translation in column B) as rows 103-138 of the translations sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEntries = @(
    @("공격 시 1% 확률로 오토스틸", "1% Auto Steal"),
    @("공격 시 2% 확률로 오토스틸", "2% Auto Steal"),
    @("공격 시 3% 확률로 오토스틸", "3% Auto Steal"),
    @("공격 시 5% 확률로 오토스틸", "5% Auto Steal"),
    @("공격 시 7% 확률로 오토스틸", "7% Auto Steal"),
    @("공격 시 3% 확률로 32의 HP 회복", "3% chance to recover 32 HP on attack"),
    @("공격 시 3% 확률로 32의 MP 회복", "3% chance to recover 32 MP on attack"),
    @("공격 시 10% 확률로 6레벨 중독효과 적용", "10% chance to apply Lv 6 Poison on attack"),
    @("공격 시 5% 확률로 2레벨 기절효과 적용", "5% chance to apply Lv 2 Stun on attack"),
    @("공격 시 10% 확률로 2레벨 슬로우효과 적용", "10% chance to apply Lv 2 Slow on attack"),
    @("공격 시 10% 확률로 3레벨 암흑효과 적용", "10% chance to apply Lv 3 Darkness on attack"),
    @("공격 시 5% 확률로 2레벨 빙결효과 적용", "5% chance to apply Lv 2 Freeze on attack"),
    @("공격 시 5% 확률로 2레벨 봉인효과 적용", "5% chance to apply Lv 2 Seal on attack"),
    @("공격 시 3% 확률로 28의 HP 회복", "3% chance to recover 28 HP on attack"),
    @("공격 시 3% 확률로 28의 MP 회복", "3% chance to recover 28 MP on attack"),
    @("공격 시 10% 확률로 5레벨 중독효과 적용", "10% chance to apply Lv 5 Poison on attack"),
    @("공격 시 5% 확률로 2레벨 기절효과 적용", "5% chance to apply Lv 2 Stun on attack"),
    @("공격 시 10% 확률로 2레벨 슬로우효과 적용", "10% chance to apply Lv 2 Slow on attack"),
    @("공격 시 10% 확률로 2레벨 암흑효과 적용", "10% chance to apply Lv 2 Darkness on attack"),
    @("공격 시 5% 확률로 2레벨 빙결효과 적용", "5% chance to apply Lv 2 Freeze on attack"),
    @("공격 시 5% 확률로 2레벨 봉인효과 적용", "5% chance to apply Lv 2 Seal on attack"),
    @("피격 시 20% 확률로 21의 데미지 무시", "20% chance to ignore 21 damage when attacked"),
    @("피격 시 20% 확률로 32의 데미지 무시", "20% chance to ignore 32 damage when attacked"),
    @("피격 시 30% 확률로 43의 데미지 무시", "30% chance to ignore 43 damage when attacked"),
    @("피격 시 20% 확률로 25의 데미지 무시", "20% chance to ignore 25 damage when attacked"),
    @("피격 시 20% 확률로 38의 데미지 무시", "20% chance to ignore 38 damage when attacked"),
    @("피격 시 30% 확률로 51의 데미지 무시", "30% chance to ignore 51 damage when attacked"),
    @("피격 시 10% 확률로 10초간 분노를 느낀다.", "10% chance to become Angry for 10 sec when attacked"),
    @("피격 시 10% 확률로 10초간 행복을 느낀다.", "10% chance to become Happy for 10 sec when attacked"),
    @("피격 시 10% 확률로 10초간 사랑에 빠진다.", "10% chance to become In Love for 10 sec when attacked"),
    @("피격 시 10% 확률로 10초간 격노를 느낀다.", "10% chance to become Enraged for 10 sec when attacked"),
    @("피격 시 10% 확률로 10초간 감동을 느낀다.", "10% chance to become Touched for 10 sec when attacked"),
    @("4초 당 24의 HP 회복", "Recover 24 HP every 4 seconds"),
    @("4초 당 24의 MP 회복", "Recover 24 MP every 4 seconds"),
    @("모든 스킬레벨", "All Skill Levels"),
    @("모든 스킬의 재사용 대기시간", "Skill Cooldown Reduction")
)

$startRow = 103
$r = $startRow
foreach ($entry in $newEntries) {
    $ws.Cells.Item($r, 1).Value = $entry[0]
    # Column B keeps the existing "Text" number format used throughout the sheet
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $r = $r + 1
}

# Restore the cursor/selection position recorded in the authored workbook
$ws.Range("F151").Select() | Out-Null
